# Auto-applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a Number by Excel
# (losing significant trailing zeros / exact decimal text). Force them to stay
# as plain text by switching the cell to a Text number format before writing.
$textCells = @("D4", "D6", "D7", "D8", "D9", "D11", "D12", "D15", "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D31", "D33", "D34", "D36", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.214.88"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.863.65"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "235.16"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4666"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "0.2831"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "0.06520"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("D11").Value = "0.07859"
$ws.Range("E11").Value = "  +1.09%  "
$ws.Range("D12").Value = "97.41"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "1.868.71"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "5.101"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "0.6726"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "280.29"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("D17").Value = "30.203.20"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "5.511"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "2.113.46"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "0.000007272"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "6.146"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "9.195"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("D26").Value = "164.79"
$ws.Range("E26").Value = "  -1.98%  "
$ws.Range("D27").Value = "19.14"
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "1.923"
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").Value = "1.379"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "4.409"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "4.089"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("D34").Value = "0.04695"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "1.115"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").Value = "0.7074"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "2.536"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").Value = "6.218"
$ws.Range("E40").Value = "  -7.32%  "
$ws.Range("D41").Value = "73.24"
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").Value = "1.943"
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("D43").Value = "0.8480"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("D44").Value = "103.99"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4164"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9998"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "7.198"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").Value = "9.133"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").Value = "933.68"
$ws.Range("E49").Value = "  -6.01%  "
$ws.Range("D50").Value = "34.17"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").Value = "0.1128"
$ws.Range("E51").Value = "  -1.83%  "
